$wb = $excel.ActiveWorkbook

# --- "Carolina Ferreira" sheet (sheet1 / rId1) -----------------------------
$ws = $wb.Worksheets.Item("Carolina Ferreira")

# Column B (heuristic counts for "Tela Principal")
$ws.Cells.Item(3,  2).Value = "-"
$ws.Cells.Item(4,  2).Value = "-"
$ws.Cells.Item(5,  2).Value = "-"
$ws.Cells.Item(6,  2).Value = 0
$ws.Cells.Item(7,  2).Value = "-"
$ws.Cells.Item(8,  2).Value = 0
$ws.Cells.Item(9,  2).Value = "-"
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(11, 2).Value = "-"
$ws.Cells.Item(12, 2).Value = "-"

# Column C (heuristic counts for "Tela Mapa")
$ws.Cells.Item(3,  3).Value = 2
$ws.Cells.Item(4,  3).Value = "-"
$ws.Cells.Item(5,  3).Value = 1
$ws.Cells.Item(6,  3).Value = 0
$ws.Cells.Item(7,  3).Value = 2
$ws.Cells.Item(8,  3).Value = "-"
$ws.Cells.Item(9,  3).Value = "-"
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(11, 3).Value = 1
$ws.Cells.Item(12, 3).Value = "-"
$ws.Cells.Item(13, 3).Value = "2 - limpar filtros"
$ws.Cells.Item(14, 3).Value = "5 - confirmar filtro"
$ws.Cells.Item(15, 3).Value = "9 - sem resultados"

# Column D (heuristic counts for "Tela Bate-Papo")
$ws.Cells.Item(3,  4).Value = "-"
$ws.Cells.Item(4,  4).Value = "-"
$ws.Cells.Item(5,  4).Value = "-"
$ws.Cells.Item(6,  4).Value = 0
$ws.Cells.Item(7,  4).Value = "-"
$ws.Cells.Item(8,  4).Value = "-"
$ws.Cells.Item(9,  4).Value = "-"
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(12, 4).Value = "-"
$ws.Cells.Item(13, 4).Value = "9 - erro de envio"

# Column E (heuristic counts for "Tela Evento")
$ws.Cells.Item(3,  5).Value = "-"
$ws.Cells.Item(4,  5).Value = "-"
$ws.Cells.Item(5,  5).Value = 1
$ws.Cells.Item(6,  5).Value = 1
$ws.Cells.Item(7,  5).Value = 2
$ws.Cells.Item(8,  5).Value = "-"
$ws.Cells.Item(9,  5).Value = "-"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(11, 5).Value = "-"
$ws.Cells.Item(13, 5).Value = "3 - desinscrever"
$ws.Cells.Item(14, 5).Value = "5 - conf. presença"

# Column F (heuristic counts for "Tela Usuário")
$ws.Cells.Item(3,  6).Value = "-"
$ws.Cells.Item(4,  6).Value = "-"
$ws.Cells.Item(5,  6).Value = "-"
$ws.Cells.Item(6,  6).Value = 1
$ws.Cells.Item(7,  6).Value = "-"
$ws.Cells.Item(8,  6).Value = "-"
$ws.Cells.Item(9,  6).Value = "-"
$ws.Cells.Item(10, 6).Value = 2
$ws.Cells.Item(11, 6).Value = "-"
$ws.Cells.Item(12, 6).Value = "-"

# Resize columns C and E to fit the new (longer) text that was typed in
$ws.Columns.Item(3).ColumnWidth = 14.666666666666666
$ws.Columns.Item(5).ColumnWidth = 14.833333333333334

# Make this the active sheet / cell, matching the saved view state (this
# also clears tabSelected on whichever sheet used to be active).
[void]$ws.Activate()
[void]$ws.Range("D18").Select()
